$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old stray "pvilje" header in A1 is gone; the participant list now
# starts on row 2 and gains a new entry ("name 16") in its proper sorted
# spot between "name 15" and "name 17", pushing the list down to end at
# row 19.
$ws.Range("A1").ClearContents()

$names = @(
    "name 1", "name 2", "name 3", "name 4", "name 5", "name 6", "name 7",
    "name 8", "name 9", "name 10", "name 11", "name 12", "name 13",
    "name 14", "name 15", "name 16", "name 17", "name 18"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $names[$i]
}

# Match the active selection recorded for the edited sheet.
$ws.Range("E6").Select()
